$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 226 ("Espinaca" weekly price records),
# which shifts the existing rows 226-268 down to 227-269.
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with this week's record.
$ws.Range("A226").Value = 8
$ws.Range("B226").Value = "Terminal La Palmera de La Serena"
$ws.Range("C226").Value = "Coquimbo"
$ws.Range("D226").Value = 44694
$ws.Range("E226").Value = 4
$ws.Range("F226").Value = 100112012
$ws.Range("G226").Value = "Espinaca"
$ws.Range("H226").Value = "Sin especificar"
$ws.Range("I226").Value = "Primera"
$ws.Range("J226").Value = 3200
$ws.Range("K226").Value = 500
$ws.Range("L226").Value = 550
$ws.Range("M226").Value = 525
$ws.Range("N226").Value = "$/atado 300 a 500 gramos"
$ws.Range("O226").Value = "Provincia del Elquí"
$ws.Range("P226").Value = 1050
$ws.Range("Q226").Value = 0.5
$ws.Range("R226").Value = "Hortaliza"
